$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team Dim")

# New "Trophies" column header
$ws.Range("N1").Value = "Trophies"

# Trophies values for each team row (rows 2-16)
$trophies = @(5, 1, 0, 0, 1, 0, 0, 3, 0, 5, 0, 1, 0, 1, 1)

for ($i = 0; $i -lt $trophies.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 14).Value = $trophies[$i]
}

# Update selection to match the state recorded in the saved workbook
$ws.Activate()
$ws.Range("N18:N19").Select()
